# Update the "想去人数" (want-to-go count) figures that changed between
# the previous gh-pages data snapshot and the newly generated one.
#
# Sheet "展览" (Exhibitions):
#   F3: 2074 -> 2081
#   F4: 858  -> 859
#   F5: 1208 -> 1228
#   F6: 357  -> 358
#
# Sheet "全部类型" (All types, aggregate of every sheet):
#   F3: 2074 -> 2081
#   F6: 858  -> 859
#   F7: 1208 -> 1228
#   F8: 357  -> 358

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item("展览")
$wsExhibitions.Range("F3").Value = 2081
$wsExhibitions.Range("F4").Value = 859
$wsExhibitions.Range("F5").Value = 1228
$wsExhibitions.Range("F6").Value = 358

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 2081
$wsAll.Range("F6").Value = 859
$wsAll.Range("F7").Value = 1228
$wsAll.Range("F8").Value = 358
